$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table data (rows 16-33), grouped by worker (instead of by period as before).
# Columns: C = Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora, F = Valor Mora, G = Salario Basico
$data = @(
    @("1143358292", "ALDAIR MAURICIO MARTINEZ MONTALVO", "2307", 36341, 1160000),
    @("1143358292", "ALDAIR MAURICIO MARTINEZ MONTALVO", "2306", 46400, 1160000),
    @("1143358292", "ALDAIR MAURICIO MARTINEZ MONTALVO", "2305", 46400, 1160000),
    @("1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL", "2308", 160000, 4000000),
    @("1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL", "2307", 46400, 1160000),
    @("1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL", "2306", 160000, 4000000),
    @("1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL", "2305", 160000, 4000000),
    @("1128050143", "LEIDYS DEL CARMEN SARA PAJARO", "2308", 80000, 2000000),
    @("1128050143", "LEIDYS DEL CARMEN SARA PAJARO", "2307", 80000, 2000000),
    @("1128050143", "LEIDYS DEL CARMEN SARA PAJARO", "2306", 80000, 2000000),
    @("1128050143", "LEIDYS DEL CARMEN SARA PAJARO", "2305", 80000, 2000000),
    @("1237439154", "YESSIKA MARIA LOPEZ YI", "2308", 46400, 1160000),
    @("1237439154", "YESSIKA MARIA LOPEZ YI", "2307", 46400, 1160000),
    @("1237439154", "YESSIKA MARIA LOPEZ YI", "2306", 46400, 1160000),
    @("1237439154", "YESSIKA MARIA LOPEZ YI", "2305", 46400, 1160000),
    @("1002428397", "AMINTA ROSA CASTRO MERCADO", "2307", 50160, 1254000),
    @("1002428397", "AMINTA ROSA CASTRO MERCADO", "2306", 50160, 1254000),
    @("1002428397", "AMINTA ROSA CASTRO MERCADO", "2305", 50160, 1254000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 3).Value = $rec[0]
    $ws.Cells.Item($row, 4).Value = $rec[1]
    $ws.Cells.Item($row, 5).Value = $rec[2]
    $ws.Cells.Item($row, 6).Value = $rec[3]
    $ws.Cells.Item($row, 7).Value = $rec[4]
}
